# Dodan okvir za tabele u excelu i napravljeni pdf
# Select the data table (A2:C18) and apply a thin box border around every
# cell (Home > Borders > All Borders), matching the manual Excel edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:C18")
[void]$rng.Select()
$rng.Borders.LineStyle = 1
